$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two LineString coordinates that changed for the plotted elements
$ws.Range("C8").Value = "LineString( 697 29, 722 29 )"
$ws.Range("C11").Value = "LineString( 696 36, 720 36 )"

# Update the active selection to match the edited cell
$ws.Range("C11").Select()
